$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.077.88"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.216.36"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.66"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.44"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.41"
$ws.Range("E10").Value = "  +7.94%  "
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.71"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0781"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.40"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.563.44"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.85"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.194.04"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.734"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "40.044.17"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.20"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.77"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.58"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.21"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.11"
$ws.Range("E28").Value = "  +2.50%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.31"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.05"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.67"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.06"
$ws.Range("E34").Value = "  +6.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.94"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0712"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0998"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.57"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.077.17"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.36"
$ws.Range("E44").Value = "  +12.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0270"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.91"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.76"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.94"
$ws.Range("E48").Value = "  -9.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.438.16"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("E51").Value = "  +1.01%  "
